$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "T GILPIN PHYSIO CONSULTANCY LTD"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "16460503"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 8).Value = "LP"

$ws.Cells.Item(3, 1).Value = "SAMVIV PARTNERS LTD"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "16460672"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 8).Value = "Partners"

$ws.Cells.Item(4, 1).Value = "4D CAPITAL PROPCO (44) LIMITED"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "16461269"
$ws.Cells.Item(4, 2).Style = "Normal"
$ws.Cells.Item(4, 8).Value = "Capital"

$ws.Cells.Item(5, 1).Value = "DGPI LTD"
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "SC849118"
$ws.Cells.Item(5, 2).Style = "Normal"
$ws.Cells.Item(5, 8).Value = "GP"

$ws.Cells.Item(6, 1).Value = "AFROSCOT VENTURES LTD"
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "16462878"
$ws.Cells.Item(6, 2).Style = "Normal"
$ws.Cells.Item(6, 8).Value = "Ventures"

$ws.Cells.Item(7, 1).Value = "ST GEORGE CAPITAL (LAND) LIMITED"
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "16462880"
$ws.Cells.Item(7, 2).Style = "Normal"
$ws.Cells.Item(7, 8).Value = "Capital"

$ws.Cells.Item(8, 1).Value = "DAVIDSON CAPITAL HOLDINGS LTD"
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "SC849117"
$ws.Cells.Item(8, 2).Style = "Normal"
$ws.Cells.Item(8, 8).Value = "Capital"

$ws.Cells.Item(9, 1).Value = "KC INVESTMENTS & TRADING LIMITED"
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "16456642"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 8).Value = "Investments"

$ws.Cells.Item(10, 1).Value = "JJOHN INVESTMENTS LIMITED"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "16456276"
$ws.Cells.Item(10, 2).Style = "Normal"
$ws.Cells.Item(10, 8).Value = "Investments"

$ws.Cells.Item(11, 1).Value = "THE REEL MED LLP"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "OC456780"
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 8).Value = "LP"

$ws.Cells.Item(12, 1).Value = "PONGPONG MALATANG LTD"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = "16458077"
$ws.Cells.Item(12, 2).Style = "Normal"
$ws.Cells.Item(12, 8).Value = "GP"

$ws.Cells.Item(13, 1).Value = "KNOTT INVESTMENTS LIMITED"
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "16458684"
$ws.Cells.Item(13, 2).Style = "Normal"
$ws.Cells.Item(13, 8).Value = "Investments"

$ws.Cells.Item(14, 1).Value = "MUSICROOTS LTD"
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "16455514"
$ws.Cells.Item(14, 2).Style = "Normal"
$ws.Cells.Item(14, 8).Value = "SIC"

$ws.Cells.Item(16, 1).Value = "ECHO VENTURES GROUP LIMITED"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "16455744"
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 8).Value = "Ventures"

$ws.Cells.Item(17, 1).Value = "TALKSGPT AI LTD"
$ws.Cells.Item(17, 2).NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = "16455313"
$ws.Cells.Item(17, 2).Style = "Normal"
$ws.Cells.Item(17, 8).Value = "GP"

$ws.Cells.Item(19, 1).Value = "DAVISON FAMILY CAPITAL LTD"
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = "16455115"
$ws.Cells.Item(19, 2).Style = "Normal"
$ws.Cells.Item(19, 8).Value = "Capital"

$ws.Cells.Item(21, 1).Value = "GROWTHFORGE MANAGEMENT LLP"
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = "OC456769"
$ws.Cells.Item(21, 2).Style = "Normal"
$ws.Cells.Item(21, 8).Value = "LP"

$ws.Cells.Item(22, 1).Value = "IX PARTNERS LLP"
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = "OC456771"
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 8).Value = "Partners"

$ws.Cells.Item(23, 1).Value = "TUERNER IMMIGRATION LLP"
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = "OC456770"
$ws.Cells.Item(23, 2).Style = "Normal"
$ws.Cells.Item(23, 8).Value = "LP"

$ws.Cells.Item(24, 1).Value = "CAPITAL & CENTRIC (SYNCHRONICITY) LTD"
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = "16453716"
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 8).Value = "Capital"

$ws.Cells.Item(25, 1).Value = "CAMBRIDGE SOCIAL INVESTMENTS LIMITED"
$ws.Cells.Item(25, 2).NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = "16453466"
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 8).Value = "Investments"

$ws.Cells.Item(26, 1).Value = "ALDABBOUS UK INVESTMENTS LTD"
$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "16453476"
$ws.Cells.Item(26, 2).Style = "Normal"
$ws.Cells.Item(26, 8).Value = "Investments"

$ws.Cells.Item(27, 1).Value = "GULF TRADE AND INVESTMENT ADVANTAGES JOINT PARTNERSHIP LTD"
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "16453733"
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 8).Value = "Partners"

$ws.Cells.Item(28, 1).Value = "GOLDEN VENTURES LONDON LTD"
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = "16452104"
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 8).Value = "Ventures"

$ws.Cells.Item(29, 1).Value = "ASSET CAPITAL 44 OPCO LIMITED"
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = "16449512"
$ws.Cells.Item(29, 2).Style = "Normal"
$ws.Cells.Item(29, 8).Value = "Capital"

$ws.Cells.Item(30, 1).Value = "SYNERGY FUNDING LTD"
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "16449538"
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 8).Value = "Fund"

$ws.Cells.Item(31, 1).Value = "ATHENA PARTNERSHIP LTD"
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "16449517"
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 8).Value = "Partners"

$ws.Cells.Item(32, 1).Value = "FROST CAPITAL LTD"
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = "16450073"
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(32, 8).Value = "Capital"
